$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EK5-PÇ Karşılama Yüzdeleri")

$cols = @(
    @("A", "test"),
    @("B", "test2"),
    @("E", 0.9248000383377075),
    @("F", 0.9248000383377075),
    @("G", 0.9248000383377075),
    @("H", 0.9248000383377075),
    @("I", 0.9248000383377075),
    @("J", 0.9371429085731506),
    @("K", 0.9248000383377075),
    @("L", 0.9248000383377075),
    @("P", 0.6489999890327454),
    @("Q", 0.6489999890327454),
    @("R", 0.6489999890327454),
    @("S", 0.6489999890327454),
    @("T", 0.6489999890327454),
    @("U", 0.6489999890327454),
    @("V", 0.6489999890327454),
    @("W", 0.6489999890327454),
    @("AA", 0.6800000071525574),
    @("AB", 0.6800000071525574),
    @("AC", 0.6800000071525574),
    @("AD", 0.6800000071525574),
    @("AE", 0.6800000071525574),
    @("AF", 0.6800000071525574),
    @("AG", 0.6800000071525574),
    @("AH", 0.6800000071525574),
    @("AJ", 0.6840000152587891),
    @("AM", 0.6840000152587891),
    @("AN", 0.6840000152587891),
    @("AO", 0.6840000152587891),
    @("AP", 0.6840000152587891),
    @("AR", 0.6840000152587891),
    @("AS", 0.6840000152587891),
    @("AU", 0.6840000152587891),
    @("AW", 0.7512667179107666),
    @("AX", 0.7344500422477722),
    @("AY", 0.7344500422477722),
    @("AZ", 0.7344500422477722),
    @("BA", 0.7344500422477722),
    @("BB", 0.7553809285163879),
    @("BC", 0.7344500422477722),
    @("BD", 0.7344500422477722)
)

for ($r = 29; $r -le 32; $r++) {
    foreach ($pair in $cols) {
        $ws.Range($pair[0] + $r).Value = $pair[1]
    }
}
